{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the index of the \"LOQ4205...\" paragraph so we can also remove the\n// blank paragraph that immediately follows it (the diff removes that blank\n// separator along with the two text paragraphs below it).\nconst items = paragraphs.items;\nlet loqIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"LOQ4205: Sistemas Produtivos II (Requisito fraco)\") {\n    loqIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (loqIndex !== -1 && loqIndex + 1 < items.length && items[loqIndex + 1].text === \"\") {\n  toDelete.push(items[loqIndex + 1]);\n}\nfor (let i = 0; i < items.length; i++) {\n  if (targetTexts.indexOf(items[i].text) !== -1) {\n    toDelete.push(items[i]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "$word = New-Object -ComObject Word.Application\n$d = $word.ActiveDocument\n\n$targetTexts = @(\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Collect the paragraphs to delete: the two \"footer\" paragraphs plus the\n# blank separator paragraph that immediately precedes the first one of them\n# (right after the \"LOQ4205...\" requirement line).\n$count = $d.Paragraphs.Count\n$toDelete = @()\n\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    $trimmed = $text.TrimEnd(\"`r\", \"`a\")\n    if ($trimmed -eq \"LOQ4205: Sistemas Produtivos II (Requisito fraco)\") {\n        if (($i + 1) -le $count) {\n            $nextText = $d.Paragraphs.Item($i + 1).Range.Text.TrimEnd(\"`r\", \"`a\")\n            if ($nextText -eq \"\") {\n                $toDelete += ($i + 1)\n            }\n        }\n    }\n    if ($targetTexts -contains $trimmed) {\n        $toDelete += $i\n    }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$toDelete = $toDelete | Sort-Object -Unique -Descending\n\nforeach ($idx in $toDelete) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
